$d = $word.ActiveDocument

# --- Title ---
$d.Content.Find.Execute(
    "Quantum Entanglement: Unveiling the Quantum Mysteries", $true, $false, $false, $false, $false,
    $true, 1, $false, "The Enchanting Realm of Biology: Unveiling the Wonders of Life", 2) | Out-Null

# --- Author ---
$d.Content.Find.Execute(
    "Thomas Langton", $true, $false, $false, $false, $false,
    $true, 1, $false, "Professor Emily Carter", 2) | Out-Null

# --- Email line: "thomas" + "." + "langton@gmail" + "." + "com"  ->  "emilycarter@hscemail" + "." + "com" ---
$d.Content.Find.Execute(
    "thomas.langton@gmail", $true, $false, $false, $false, $false,
    $true, 1, $false, "emilycarter@hscemail", 2) | Out-Null

# --- Body paragraph (intro) ---
$d.Content.Find.Execute(
    "In the realm of quantum physics, the enigmatic phenomenon known as quantum entanglement defies our classical understanding of reality",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Biology, the study of life, embarks on a thrilling voyage into the captivating realm of living organisms", 2) | Out-Null

$d.Content.Find.Execute(
    " This bizarre and counterintuitive correlation between particles separated by vast distances has captivated the minds of scientists for decades",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " It unveils the intricate intricacies of cells, uncovers the mysteries embedded within DNA, and delves into the extraordinary diversity of life forms inhabiting our planet", 2) | Out-Null

$d.Content.Find.Execute(
    " As we delve into the intricate tapestry of quantum entanglement, we embark on a journey to unravel its profound implications, exploring its potential to revolutionize fields like quantum computing, cryptography, and even our understanding of the universe itself",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " From the microscopic world of bacteria to the majestic grandeur of whales, biology unlocks the secrets of existence, revealing the profound interconnectedness that binds all living things. In this exploration of life's wonders, we embark on an exhilarating journey, unraveling the tapestry of life's rich symphony", 2) | Out-Null

$d.Content.Find.Execute(
    "The essence of quantum entanglement lies in the interconnectedness of particles, such that the state of one particle instantaneously affects the state of its entangled partner, regardless of the distance separating them",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "The study of genetics unveils the intricate mechanisms of heredity, deciphering the genetic blueprints that govern the traits and characteristics passed down through generations", 2) | Out-Null

$d.Content.Find.Execute(
    " This non-local connection challenges our notion of locality, suggesting that information can travel faster than the speed of light",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " We delve into the realm of evolution, tracing the grand narrative of life's transformation across eons, driven by the forces of natural selection", 2) | Out-Null

$d.Content.Find.Execute(
    " Physicists have conducted numerous experiments that have repeatedly confirmed the existence of quantum entanglement, solidifying its status as a cornerstone of quantum theory",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " Through the lens of ecology, we unravel the delicate balance of ecosystems, revealing the intricate web of interactions that sustain the harmony of life", 2) | Out-Null

$d.Content.Find.Execute(
    "The potential applications of quantum entanglement are vast",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Biology, however, is not merely an intellectual pursuit; it holds profound implications for our daily lives and the future of our planet", 2) | Out-Null

$d.Content.Find.Execute(
    " It holds the promise of revolutionizing computation through the development of quantum computers, which harness the unique properties of entangled particles to perform calculations exponentially faster than classical computers",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " It empowers us to understand the human body, enabling us to devise innovative treatments for diseases and ameliorate human suffering", 2) | Out-Null

$d.Content.Find.Execute(
    " Quantum entanglement also offers the prospect of unbreakable cryptography, as eavesdropping on an entangled communication channel would disrupt the entanglement and thus alert the sender and receiver to the presence of an unauthorized third party",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " By comprehending the intricate workings of ecosystems, we can devise strategies to protect and preserve the delicate balance of nature. As we continue to unravel the mysteries of life, biology empowers us to confront global challenges such as food security, climate change, and the preservation of biodiversity", 2) | Out-Null

# --- Summary paragraph ---
$d.Content.Find.Execute(
    "Quantum entanglement, a perplexing phenomenon in the domain of quantum physics, challenges our classical intuitions with its non-local correlations between particles",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Biology, the study of life, unveils the intricacies of living organisms, delving into the microscopic world of cells, deciphering the genetic blueprints of DNA, and exploring the magnificent diversity of life forms", 2) | Out-Null

$d.Content.Find.Execute(
    " The profound implications of entanglement stretch far beyond the realm of theoretical physics, extending to potential applications in quantum computing, cryptography, and even our comprehension of the universe",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " By unraveling the tapestry of life's symphony, biology empowers us to understand the mechanisms of heredity, trace the narrative of evolution, and unravel the delicate balance of ecosystems", 2) | Out-Null

$d.Content.Find.Execute(
    " As we continue to probe the depths of this enigmatic phenomenon, we stand at the threshold of a new era, poised to unravel the mysteries of quantum entanglement and its transformative power",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " It holds profound implications for our daily lives and the future of our planet, enabling us to devise treatments for diseases, protect the environment, and confront global challenges. Biology's journey of discovery continues to inspire awe and wonder, revealing the interconnectedness of all living things and the profound beauty of the natural world", 2) | Out-Null

# --- New trailing empty paragraph at the end of the document body ---
$endRange = $d.Content
$endRange.Collapse(0) | Out-Null
$endRange.InsertParagraphAfter() | Out-Null
